$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1640.8235
$ws.Range("I106").Value = 1354.4615
$ws.Range("K106").Value = 1354.4615
$ws.Range("M106").Value = -723.4614999999999

$ws.Range("H132").Value = 33861.375
$ws.Range("I132").Value = 41333.617
$ws.Range("J132").Value = 1481.6666
$ws.Range("K132").Value = 124000.851
$ws.Range("L132").Value = 4444.9998
$ws.Range("M132").Value = -121470.851
$ws.Range("N132").Value = -9504.9998

$ws.Range("H137").Value = 17172.691
$ws.Range("I137").Value = 1980.2703
$ws.Range("K137").Value = 5940.810899999999
$ws.Range("M137").Value = -3390.810899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12913.122
$ws.Range("I32").Value = 13544.72
$ws.Range("J32").Value = 6439.25
$ws.Range("K32").Value = 13544.72
$ws.Range("L32").Value = 6439.25
$ws.Range("M32").Value = -13257.72
$ws.Range("N32").Value = -7013.25

$ws.Range("H61").Value = 696863.25
$ws.Range("I61").Value = 1804811.8
$ws.Range("J61").Value = 4395.4375
$ws.Range("K61").Value = 1804811.8
$ws.Range("L61").Value = 4395.4375
$ws.Range("M61").Value = -1804599.8
$ws.Range("N61").Value = -4819.4375

$ws.Range("H74").Value = 3190.25
$ws.Range("I74").Value = 5612.375
$ws.Range("J74").Value = 1575.5
$ws.Range("K74").Value = 5612.375
$ws.Range("L74").Value = 1575.5
$ws.Range("M74").Value = -4738.375
$ws.Range("N74").Value = -3323.5

$ws.Range("H77").Value = 3190.25
$ws.Range("I77").Value = 5612.375
$ws.Range("J77").Value = 1575.5
$ws.Range("K77").Value = 28061.875
$ws.Range("L77").Value = 7877.5
$ws.Range("M77").Value = -23693.875
$ws.Range("N77").Value = -16613.5

$ws.Range("H97").Value = 2217.2856
$ws.Range("J97").Value = 905.5
$ws.Range("L97").Value = 905.5
$ws.Range("N97").Value = -1897.5

$ws.Range("H110").Value = 1601.9
$ws.Range("I110").Value = 1302.5333
$ws.Range("K110").Value = 1302.5333
$ws.Range("M110").Value = 742.4666999999999

$ws.Range("H122").Value = 2028.8269
$ws.Range("I122").Value = 1912.683
$ws.Range("J122").Value = 2461.7273
$ws.Range("K122").Value = 5738.049
$ws.Range("L122").Value = 7385.1819
$ws.Range("M122").Value = -3288.049
$ws.Range("N122").Value = -12285.1819

$ws.Range("H132").Value = 25093.637
$ws.Range("I132").Value = 1853.4286
$ws.Range("J132").Value = 65764
$ws.Range("K132").Value = 5560.2858
$ws.Range("L132").Value = 197292
$ws.Range("M132").Value = -3030.2858
$ws.Range("N132").Value = -202352

$ws.Range("H136").Value = 696863.25
$ws.Range("I136").Value = 1804811.8
$ws.Range("J136").Value = 4395.4375
$ws.Range("K136").Value = 5414435.4
$ws.Range("L136").Value = 13186.3125
$ws.Range("M136").Value = -5411885.4
$ws.Range("N136").Value = -18286.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1703.76
$ws.Range("I86").Value = 1520
$ws.Range("J86").Value = 1979.4
$ws.Range("K86").Value = 1520
$ws.Range("L86").Value = 1979.4
$ws.Range("M86").Value = -397
$ws.Range("N86").Value = -4225.4

$ws.Range("H89").Value = 1703.76
$ws.Range("I89").Value = 1520
$ws.Range("J89").Value = 1979.4
$ws.Range("K89").Value = 7600
$ws.Range("L89").Value = 9897
$ws.Range("M89").Value = -1984
$ws.Range("N89").Value = -21129

$ws.Range("H94").Value = 3416.913
$ws.Range("I94").Value = 1365.75
$ws.Range("K94").Value = 1365.75
$ws.Range("M94").Value = -914.75

$ws.Range("H105").Value = 2084872.1
$ws.Range("I105").Value = 1387.8572
$ws.Range("J105").Value = 5001750
$ws.Range("K105").Value = 1387.8572
$ws.Range("L105").Value = 5001750
$ws.Range("M105").Value = 359.1428000000001
$ws.Range("N105").Value = -5005244

$ws.Range("H107").Value = 616
$ws.Range("I107").Value = 600.9167
$ws.Range("J107").Value = 706.5
$ws.Range("K107").Value = 600.9167
$ws.Range("L107").Value = 706.5
$ws.Range("M107").Value = 1319.0833
$ws.Range("N107").Value = -4546.5

$ws.Range("H134").Value = 35082.934
$ws.Range("I134").Value = 41617.58
$ws.Range("J134").Value = 1102.8
$ws.Range("K134").Value = 124852.74
$ws.Range("L134").Value = 3308.4
$ws.Range("M134").Value = -122317.74
$ws.Range("N134").Value = -8378.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10609.915
$ws.Range("I31").Value = 22066.525
$ws.Range("J31").Value = 2835.7856
$ws.Range("K31").Value = 22066.525
$ws.Range("L31").Value = 2835.7856
$ws.Range("M31").Value = -21771.525
$ws.Range("N31").Value = -3425.7856

$ws.Range("H34").Value = 10609.915
$ws.Range("I34").Value = 22066.525
$ws.Range("J34").Value = 2835.7856
$ws.Range("K34").Value = 22066.525
$ws.Range("L34").Value = 2835.7856
$ws.Range("M34").Value = -21864.525
$ws.Range("N34").Value = -3239.7856

$ws.Range("H58").Value = 21251.52
$ws.Range("I58").Value = 1359.5
$ws.Range("J58").Value = 72402.42999999999
$ws.Range("K58").Value = 1359.5
$ws.Range("L58").Value = 72402.42999999999
$ws.Range("M58").Value = -1156.5
$ws.Range("N58").Value = -72808.42999999999

$ws.Range("H86").Value = 11550.25
$ws.Range("I86").Value = 5422.846
$ws.Range("K86").Value = 5422.846
$ws.Range("M86").Value = -4299.846

$ws.Range("H89").Value = 11550.25
$ws.Range("I89").Value = 5422.846
$ws.Range("K89").Value = 27114.23
$ws.Range("M89").Value = -21498.23

$ws.Range("H105").Value = 12502996
$ws.Range("J105").Value = 4500
$ws.Range("L105").Value = 4500
$ws.Range("N105").Value = -7994

$ws.Range("H132").Value = 14359.81
$ws.Range("I132").Value = 17715.807
$ws.Range("J132").Value = 4902
$ws.Range("K132").Value = 53147.421
$ws.Range("L132").Value = 14706
$ws.Range("M132").Value = -50617.421
$ws.Range("N132").Value = -19766

$ws.Range("H134").Value = 608.9524
$ws.Range("I134").Value = 516.6316
$ws.Range("K134").Value = 1549.8948
$ws.Range("M134").Value = 985.1052

$ws.Range("H136").Value = 21251.52
$ws.Range("I136").Value = 1359.5
$ws.Range("J136").Value = 72402.42999999999
$ws.Range("K136").Value = 4078.5
$ws.Range("L136").Value = 217207.29
$ws.Range("M136").Value = -1528.5
$ws.Range("N136").Value = -222307.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3588
$ws.Range("I3").Value = 1383.3334
$ws.Range("J3").Value = 7997.3335
$ws.Range("K3").Value = 4150.0002
$ws.Range("L3").Value = 23992.0005
$ws.Range("M3").Value = -4038.0002
$ws.Range("N3").Value = -24216.0005

$ws.Range("H4").Value = 4285789.5
$ws.Range("J4").Value = 10000033
$ws.Range("L4").Value = 30000099
$ws.Range("N4").Value = -30000323

$ws.Range("H68").Value = 4243.7095
$ws.Range("J68").Value = 5088.28
$ws.Range("L68").Value = 15264.84
$ws.Range("N68").Value = -16886.84

$ws.Range("H71").Value = 4243.7095
$ws.Range("J71").Value = 5088.28
$ws.Range("L71").Value = 45794.52
$ws.Range("N71").Value = -53906.52

$ws.Range("H104").Value = 36676.332
$ws.Range("J104").Value = 36676.332
$ws.Range("L104").Value = 110028.996
$ws.Range("N104").Value = -115270.996

$ws.Range("H131").Value = 799.26
$ws.Range("J131").Value = 807
$ws.Range("L131").Value = 2421
$ws.Range("N131").Value = -12501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8040.6665
$ws.Range("I80").Value = 15213.125
$ws.Range("J80").Value = 3626.8462
$ws.Range("K80").Value = 15213.125
$ws.Range("L80").Value = 3626.8462
$ws.Range("M80").Value = -14215.125
$ws.Range("N80").Value = -5622.8462

$ws.Range("H83").Value = 8040.6665
$ws.Range("I83").Value = 15213.125
$ws.Range("J83").Value = 3626.8462
$ws.Range("K83").Value = 76065.625
$ws.Range("L83").Value = 18134.231
$ws.Range("M83").Value = -71073.625
$ws.Range("N83").Value = -28118.231

$ws.Range("H92").Value = 9217.666999999999
$ws.Range("J92").Value = 9217.666999999999
$ws.Range("L92").Value = 9217.666999999999
$ws.Range("N92").Value = -12961.667

$ws.Range("H102").Value = 1672
$ws.Range("I102").Value = 1590.3
$ws.Range("K102").Value = 1590.3
$ws.Range("M102").Value = 31.70000000000005

$ws.Range("H122").Value = 2344.652
$ws.Range("I122").Value = 2189.9443
$ws.Range("J122").Value = 2901.6
$ws.Range("K122").Value = 6569.8329
$ws.Range("L122").Value = 8704.799999999999
$ws.Range("M122").Value = -4119.8329
$ws.Range("N122").Value = -13604.8

$ws.Range("H132").Value = 57053.465
$ws.Range("I132").Value = 66316.31
$ws.Range("J132").Value = 44703
$ws.Range("K132").Value = 198948.93
$ws.Range("L132").Value = 134109
$ws.Range("M132").Value = -196418.93
$ws.Range("N132").Value = -139169

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3290.4614
$ws.Range("J68").Value = 3873.5293
$ws.Range("L68").Value = 3873.5293
$ws.Range("N68").Value = -5371.5293

$ws.Range("H71").Value = 3290.4614
$ws.Range("J71").Value = 3873.5293
$ws.Range("L71").Value = 19367.6465
$ws.Range("N71").Value = -26855.6465

$ws.Range("H132").Value = 1625.7441
$ws.Range("I132").Value = 1065.7812
$ws.Range("J132").Value = 3254.7273
$ws.Range("K132").Value = 3197.3436
$ws.Range("L132").Value = 9764.1819
$ws.Range("M132").Value = -667.3435999999997
$ws.Range("N132").Value = -14824.1819

$ws.Range("H136").Value = 19891.822
$ws.Range("I136").Value = 29137.055
$ws.Range("J136").Value = 3250.4
$ws.Range("K136").Value = 87411.16500000001
$ws.Range("L136").Value = 9751.200000000001
$ws.Range("M136").Value = -84861.16500000001
$ws.Range("N136").Value = -14851.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 9602
$ws.Range("J103").Value = 9602
$ws.Range("L103").Value = 9602
$ws.Range("N103").Value = -11946

$ws.Range("H107").Value = 1934.5
$ws.Range("I107").Value = 639
$ws.Range("J107").Value = 2452.7
$ws.Range("K107").Value = 1917
$ws.Range("L107").Value = 7358.099999999999
$ws.Range("M107").Value = 3
$ws.Range("N107").Value = -11198.1

$ws.Range("H132").Value = 1834.4117
$ws.Range("I132").Value = 1598.0555
$ws.Range("J132").Value = 2401.6667
$ws.Range("K132").Value = 4794.166499999999
$ws.Range("L132").Value = 7205.000100000001
$ws.Range("M132").Value = -2264.166499999999
$ws.Range("N132").Value = -12265.0001
